$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.705.14'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '3.384.02'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.05%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.384.35'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.475'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.400'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').Value = '3.962.60'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '3.381.35'
$ws.Range('E17').Value = '  -2.38%  '
$ws.Range('D18').Value = '60.815.96'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('E21').Value = '  -5.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '384.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('E26').Value = '  -4.73%  '
$ws.Range('D27').Value = '3.521.83'
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('E28').Value = '  -1.98%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.93%  '
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.70'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '166.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.35%  '
$ws.Range('B39').Value = 'RenzoRestakedETH'
$ws.Range('C39').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D39').Value = '3.415.30'
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.81'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.56%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('D48').Value = '2.520.63'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('E49').Value = '  -3.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.27%  '
